$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-02 Sunday" "2025-03-03 Monday"

Replace-Text "991÷6=" "588÷9="
Replace-Text "400÷2=" "189÷9="
Replace-Text "881÷5=" "317÷5="
Replace-Text "796÷7=" "979÷6="
Replace-Text "508÷8=" "312÷2="

Replace-Text "385÷2=" "602÷7="
Replace-Text "501÷4=" "846÷4="
Replace-Text "198÷2=" "131÷4="
Replace-Text "830÷4=" "633÷4="
Replace-Text "698÷8=" "331÷9="

Replace-Text "800÷4=" "556÷4="
Replace-Text "840÷6=" "332÷3="
Replace-Text "830÷6=" "466÷6="
Replace-Text "861÷6=" "764÷7="
Replace-Text "354÷8=" "558÷2="

Replace-Text "204÷9=" "325÷6="
Replace-Text "354÷6=" "321÷8="
Replace-Text "345÷2=" "695÷3="
Replace-Text "690÷5=" "991÷7="
Replace-Text "835÷8=" "728÷8="

Replace-Text "286÷5=" "639÷6="
Replace-Text "620÷3=" "858÷4="
Replace-Text "672÷9=" "260÷2="
Replace-Text "854÷7=" "281÷8="
Replace-Text "863÷3=" "156÷7="

Write-Output "Done replacing text"
